$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 22, shifting existing rows 22-101 down to 23-102
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the latest weekly price entry
$ws.Cells.Item(22, 1).Value = 9
$ws.Cells.Item(22, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(22, 3).Value = "Metropolitana"
$ws.Cells.Item(22, 4).Value = 45063
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 13
$ws.Cells.Item(22, 6).Value = 100112029
$ws.Cells.Item(22, 7).Value = "Orégano"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 16
$ws.Cells.Item(22, 11).Value = 18000
$ws.Cells.Item(22, 12).Value = 18000
$ws.Cells.Item(22, 13).Value = 18000
$ws.Cells.Item(22, 14).Value = "$/docena de atados"
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 6000
$ws.Cells.Item(22, 17).Value = 3
$ws.Cells.Item(22, 18).Value = "Hortaliza"
